$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number need an explicit
# text format, otherwise Excel auto-converts the assigned string into a
# numeric value (losing the original string formatting, e.g. "142.00").
$ws.Range("D2").Value = "60.948.84"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "3.383.38"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.16"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.00"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.62"
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("D12").Value = "3.963.49"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.81"
$ws.Range("E14").Value = "  -1.92%  "
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "3.366.54"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").Value = "61.044.41"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("E18").Value = "  -3.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.65"
$ws.Range("E19").Value = "  -4.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.98"
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.21"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.20"
$ws.Range("E22").Value = "  +2.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.551"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -3.00%  "
$ws.Range("D26").Value = "3.522.38"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.27"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.00"
$ws.Range("E30").Value = "  -2.13%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -4.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.30"
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.96"
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.68"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").Value = "3.414.52"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("E40").Value = "  -1.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.70"
$ws.Range("E41").Value = "  -1.65%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("E45").Value = "  -1.78%  "
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("D47").Value = "2.458.18"
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.72"
$ws.Range("E49").Value = "  -2.60%  "
$ws.Range("E50").Value = "  +8.56%  "
